# Auto-generated edit script applying the diff to khl_referees_stats_1369.xlsx
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Главные")
$ws.Range("C2").Value = 32
$ws.Range("D2").Value = 689
$ws.Range("E2").Value = 291
$ws.Range("F2").Value = 398
$ws.Range("G2").Value = 21.53
$ws.Range("H2").Value = 9.09
$ws.Range("I2").Value = 12.44
$ws.Range("J2").Value = 128
$ws.Range("K2").Value = 154
$ws.Range("V2").Value = 12
$ws.Range("AA2").Value = '2025-12-01 03:05:18'
$ws.Range("AA3").Value = '2025-12-01 03:05:18'
$ws.Range("AA4").Value = '2025-12-01 03:05:18'
$ws.Range("AA5").Value = '2025-12-01 03:05:18'
$ws.Range("AA6").Value = '2025-12-01 03:05:18'
$ws.Range("AA7").Value = '2025-12-01 03:05:18'
$ws.Range("C8").Value = 27
$ws.Range("D8").Value = 496
$ws.Range("E8").Value = 263
$ws.Range("F8").Value = 233
$ws.Range("G8").Value = 18.37
$ws.Range("H8").Value = 9.74
$ws.Range("I8").Value = 8.630000000000001
$ws.Range("J8").Value = 114
$ws.Range("K8").Value = 109
$ws.Range("P8").Value = 1
$ws.Range("AA8").Value = '2025-12-01 03:05:18'
$ws.Range("C9").Value = 30
$ws.Range("D9").Value = 480
$ws.Range("E9").Value = 249
$ws.Range("F9").Value = 231
$ws.Range("G9").Value = 16
$ws.Range("H9").Value = 8.300000000000001
$ws.Range("I9").Value = 7.7
$ws.Range("J9").Value = 122
$ws.Range("K9").Value = 113
$ws.Range("V9").Value = 18
$ws.Range("AA9").Value = '2025-12-01 03:05:18'
$ws.Range("AA10").Value = '2025-12-01 03:05:18'
$ws.Range("C11").Value = 22
$ws.Range("D11").Value = 532
$ws.Range("E11").Value = 244
$ws.Range("F11").Value = 288
$ws.Range("G11").Value = 24.18
$ws.Range("H11").Value = 11.09
$ws.Range("I11").Value = 13.09
$ws.Range("J11").Value = 107
$ws.Range("K11").Value = 99
$ws.Range("AA11").Value = '2025-12-01 03:05:18'
$ws.Range("AA12").Value = '2025-12-01 03:05:18'
$ws.Range("AA13").Value = '2025-12-01 03:05:18'
$ws.Range("AA14").Value = '2025-12-01 03:05:18'
$ws.Range("AA15").Value = '2025-12-01 03:05:18'
$ws.Range("AA16").Value = '2025-12-01 03:05:18'
$ws.Range("AA17").Value = '2025-12-01 03:05:18'
$ws.Range("C18").Value = 28
$ws.Range("D18").Value = 450
$ws.Range("E18").Value = 225
$ws.Range("F18").Value = 225
$ws.Range("G18").Value = 16.07
$ws.Range("H18").Value = 8.039999999999999
$ws.Range("I18").Value = 8.039999999999999
$ws.Range("J18").Value = 90
$ws.Range("K18").Value = 105
$ws.Range("P18").Value = 1
$ws.Range("AA18").Value = '2025-12-01 03:05:18'
$ws.Range("AA19").Value = '2025-12-01 03:05:18'
$ws.Range("AA20").Value = '2025-12-01 03:05:18'
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 362
$ws.Range("E21").Value = 160
$ws.Range("F21").Value = 202
$ws.Range("G21").Value = 13.92
$ws.Range("H21").Value = 6.15
$ws.Range("I21").Value = 7.77
$ws.Range("J21").Value = 70
$ws.Range("K21").Value = 86
$ws.Range("AA21").Value = '2025-12-01 03:05:18'
$ws.Range("AA22").Value = '2025-12-01 03:05:18'
$ws.Range("AA23").Value = '2025-12-01 03:05:18'
$ws.Range("AA24").Value = '2025-12-01 03:05:18'
$ws.Range("AA25").Value = '2025-12-01 03:05:18'
$ws.Range("AA26").Value = '2025-12-01 03:05:18'

$ws = $wb.Worksheets.Item("Линейные")
$ws.Range("AA2").Value = '2025-12-01 03:05:18'
$ws.Range("C3").Value = 29
$ws.Range("D3").Value = 443
$ws.Range("E3").Value = 241
$ws.Range("F3").Value = 202
$ws.Range("G3").Value = 15.28
$ws.Range("H3").Value = 8.31
$ws.Range("I3").Value = 6.97
$ws.Range("J3").Value = 108
$ws.Range("K3").Value = 86
$ws.Range("P3").Value = 1
$ws.Range("AA3").Value = '2025-12-01 03:05:18'
$ws.Range("AA4").Value = '2025-12-01 03:05:18'
$ws.Range("AA5").Value = '2025-12-01 03:05:18'
$ws.Range("AA6").Value = '2025-12-01 03:05:18'
$ws.Range("AA7").Value = '2025-12-01 03:05:18'
$ws.Range("AA8").Value = '2025-12-01 03:05:18'
$ws.Range("AA9").Value = '2025-12-01 03:05:18'
$ws.Range("AA10").Value = '2025-12-01 03:05:18'
$ws.Range("AA11").Value = '2025-12-01 03:05:18'
$ws.Range("C12").Value = 25
$ws.Range("D12").Value = 438
$ws.Range("E12").Value = 212
$ws.Range("F12").Value = 226
$ws.Range("G12").Value = 17.52
$ws.Range("H12").Value = 8.48
$ws.Range("I12").Value = 9.039999999999999
$ws.Range("J12").Value = 96
$ws.Range("K12").Value = 103
$ws.Range("AA12").Value = '2025-12-01 03:05:18'
$ws.Range("AA13").Value = '2025-12-01 03:05:18'
$ws.Range("AA14").Value = '2025-12-01 03:05:18'
$ws.Range("C15").Value = 25
$ws.Range("D15").Value = 489
$ws.Range("E15").Value = 251
$ws.Range("F15").Value = 238
$ws.Range("G15").Value = 19.56
$ws.Range("H15").Value = 10.04
$ws.Range("I15").Value = 9.52
$ws.Range("J15").Value = 103
$ws.Range("K15").Value = 99
$ws.Range("V15").Value = 14
$ws.Range("AA15").Value = '2025-12-01 03:05:18'
$ws.Range("AA16").Value = '2025-12-01 03:05:18'
$ws.Range("AA17").Value = '2025-12-01 03:05:18'
$ws.Range("AA18").Value = '2025-12-01 03:05:18'
$ws.Range("AA19").Value = '2025-12-01 03:05:18'
$ws.Range("AA20").Value = '2025-12-01 03:05:18'
$ws.Range("AA21").Value = '2025-12-01 03:05:18'
$ws.Range("AA22").Value = '2025-12-01 03:05:18'
$ws.Range("AA23").Value = '2025-12-01 03:05:18'
$ws.Range("AA24").Value = '2025-12-01 03:05:18'
$ws.Range("AA25").Value = '2025-12-01 03:05:18'
$ws.Range("AA26").Value = '2025-12-01 03:05:18'

Write-Output "Edit complete."